# Update the "Notes" section of the triangle-checker overview sheet with the
# algorithm for all the triangle conditions (matches commit:
# "update the js file with algorithm for all conditions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 22: duplicate of the "triangle inequality" condition (same text already
# used lower down in A28) placed near the isosceles-condition block.
$ws.Range("A22").Value = "userInput1+userInput2 >userInput3 && userInput1+userInput3 >userInput2 && userInput3+userInput2 >userInput1"

# Rows 36-37: the parenthesised version of the triangle-inequality condition,
# written out across two lines. Populate A37 before A36 so that the new
# shared-string table entries come out in the same order as the workbook
# being reproduced ( ")" before the "(userInput1..." text ).
$ws.Range("A37").Value = ")"
$ws.Range("A36").Value = "(userInput1+userInput2 >userInput3) && (userInput1+userInput3 >userInput2) && (userInput3+userInput2 >userInput1"

# Row 33: new bold section header introducing the triangle-inequality check,
# styled the same way as the other bold headers in the sheet (isosceles /
# scalene / equilateral).
$ws.Range("A33").Value = "check on whether lengths make a triangle"
$ws.Range("A33").Font.Bold = $true

# Reflect the new content in the window's scroll position / selection, same
# as the author's saved view (top-left around row 28, active cell A33).
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A33").Select()
